$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 270, which shifts the existing rows
# 270..369 down to 271..370 (dates/values in column D/J/K/L/M/N/P/etc.
# all move down by one row, matching the rest of the table).
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row 270 with the new "Ajo" record.
$ws.Cells.Item(270, 1).Value = 4
$ws.Cells.Item(270, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(270, 3).Value = "Los Lagos"
$ws.Cells.Item(270, 4).Value = 44900
$ws.Cells.Item(270, 5).Value = 10
$ws.Cells.Item(270, 6).Value = 100112003
$ws.Cells.Item(270, 7).Value = "Ajo"
$ws.Cells.Item(270, 8).Value = "Chino"
$ws.Cells.Item(270, 9).Value = "Primera"
$ws.Cells.Item(270, 10).Value = 100
$ws.Cells.Item(270, 11).Value = 18000
$ws.Cells.Item(270, 12).Value = 20000
$ws.Cells.Item(270, 13).Value = 19000
$ws.Cells.Item(270, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(270, 15).Value = "China"
$ws.Cells.Item(270, 16).Value = 1900
$ws.Cells.Item(270, 17).Value = 10
$ws.Cells.Item(270, 18).Value = "Hortaliza"
